# Applies the crypto price/volume refresh described in the commit diff.
# Only the cells that actually changed are touched; column D (Price) values
# are forced to remain text (matching the source inlineStr cells) by using
# an apostrophe text-prefix, then the style is restored to Normal so no
# stray number formatting is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'" + '66.892.84'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.97%  '
# Row 3
$ws.Range("D3").Value = "'" + '3.106.47'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +5.18%  '
# Row 4
$ws.Range("E4").Value = '  +0.04%  '
# Row 5
$ws.Range("D5").Value = "'" + '580.71'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.58%  '
# Row 6
$ws.Range("D6").Value = "'" + '172.62'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.42%  '
# Row 7
$ws.Range("E7").Value = '  +0.02%  '
# Row 8
$ws.Range("D8").Value = "'" + '3.101.75'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +5.18%  '
# Row 9
$ws.Range("E9").Value = '  +1.43%  '
# Row 10
$ws.Range("D10").Value = "'" + '6.53'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.62%  '
# Row 11
$ws.Range("E11").Value = '  +3.59%  '
# Row 12
$ws.Range("D12").Value = "'" + '0.484'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.92%  '
# Row 13
$ws.Range("E13").Value = '  +2.06%  '
# Row 14
$ws.Range("D14").Value = "'" + '37.38'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +7.54%  '
# Row 15
$ws.Range("E15").Value = '  +0.06%  '
# Row 16
$ws.Range("D16").Value = "'" + '3.622.97'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.22%  '
# Row 17
$ws.Range("D17").Value = "'" + '66.895.43'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.95%  '
# Row 18
$ws.Range("D18").Value = "'" + '7.19'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.39%  '
# Row 19
$ws.Range("D19").Value = "'" + '3.109.34'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.43%  '
# Row 20
$ws.Range("D20").Value = "'" + '16.08'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.78%  '
# Row 21
$ws.Range("D21").Value = "'" + '485.03'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +8.86%  '
# Row 22
$ws.Range("D22").Value = "'" + '0.715'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.91%  '
# Row 23
$ws.Range("D23").Value = "'" + '7.52'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.09%  '
# Row 24
$ws.Range("D24").Value = "'" + '84.17'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.59%  '
# Row 25
$ws.Range("D25").Value = "'" + '2.39'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +6.38%  '
# Row 26
$ws.Range("D26").Value = "'" + '13.14'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +6.57%  '
# Row 27
$ws.Range("D27").Value = "'" + '10.09'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.22%  '
# Row 28
$ws.Range("E28").Value = '  -0.05%  '
# Row 29
$ws.Range("D29").Value = "'" + '7.98'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.27%  '
# Row 30
$ws.Range("D30").Value = "'" + '2.40'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.26%  '
# Row 31
$ws.Range("D31").Value = "'" + '2.69'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.69%  '
# Row 32
$ws.Range("B32").Value = 'PEPE'
$ws.Range("C32").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D32").Value = "'" + '0.0000102'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.66%  '
# Row 33
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").Value = "'" + '29.00'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +6.61%  '
# Row 34
$ws.Range("E34").Value = '  +1.40%  '
# Row 35
$ws.Range("D35").Value = "'" + '1.00'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.05%  '
# Row 36
$ws.Range("D36").Value = "'" + '5.92'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.40%  '
# Row 37
$ws.Range("E37").Value = '  +2.74%  '
# Row 38
$ws.Range("B38").Value = 'Arweave'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D38").Value = "'" + '47.92'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.15%  '
# Row 39
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").Value = "'" + '2.13'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +6.58%  '
# Row 40
$ws.Range("D40").Value = "'" + '0.318'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.01%  '
# Row 41
$ws.Range("D41").Value = "'" + '50.18'
$ws.Range("D41").Style = "Normal"
# Row 42
$ws.Range("E42").Value = '  -0.19%  '
# Row 43
$ws.Range("D43").Value = "'" + '8.68'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.68%  '
# Row 44
$ws.Range("D44").Value = "'" + '2.81'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.50%  '
# Row 45
$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").Value = "'" + '0.0362'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.86%  '
# Row 46
$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").Value = "'" + '2.839.42'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.95%  '
# Row 47
$ws.Range("D47").Value = "'" + '383.83'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.29%  '
# Row 48
$ws.Range("D48").Value = "'" + '134.95'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.36%  '
# Row 50
$ws.Range("D50").Value = "'" + '24.98'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.59%  '
# Row 51
$ws.Range("E51").Value = '  +2.87%  '
